$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 189 (pushes existing rows 189-216 down to 190-217)
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row 189 with the new weekly data point
$ws.Cells.Item(189, 1).Value = 8
$ws.Cells.Item(189, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(189, 3).Value = "Coquimbo"
$ws.Cells.Item(189, 4).Value = 45180
$ws.Cells.Item(189, 5).Value = 4
$ws.Cells.Item(189, 6).Value = 100112052
$ws.Cells.Item(189, 7).Value = "Albahaca"
$ws.Cells.Item(189, 8).Value = "Sin especificar"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 900
$ws.Cells.Item(189, 11).Value = 3500
$ws.Cells.Item(189, 12).Value = 4000
$ws.Cells.Item(189, 13).Value = 3750
$ws.Cells.Item(189, 14).Value = "$/paquete"
$ws.Cells.Item(189, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(189, 16).Value = 3750
$ws.Cells.Item(189, 17).Value = 1
$ws.Cells.Item(189, 18).Value = "Hortaliza"
